{"js": "// The document currently has no explicit styles part (word/styles.xml is\n// absent) even though its single paragraph implicitly uses the built-in\n// \"Normal\" paragraph style. This change makes that style explicit by\n// registering it on the document's style collection, which causes a\n// word/styles.xml part (containing a \"Normal\" paragraph style definition)\n// to be written out on save. No existing paragraph/run content is touched.\ncontext.document.addStyle(\"Normal\", Word.StyleType.paragraph);\nawait context.sync();\n", "ps1": "# The document currently has no explicit styles part (word/styles.xml is\n# absent) even though its single paragraph implicitly uses the built-in\n# \"Normal\" paragraph style. This change makes that style explicit by\n# registering it in the document's Styles collection, which causes a\n# word/styles.xml part (containing a \"Normal\" paragraph style definition)\n# to be written out on save. No existing paragraph/run content is touched.\n$d = $word.ActiveDocument\n# 1 == wdStyleTypeParagraph\n$d.Styles.Add(\"Normal\", 1) | Out-Null\n"}
